# Update the Belgium RES mapping scheme: reduce the share of "W" (wood)
# construction for the two "Single-family" + related columns on both
# sheets, following local expert feedback.

$wb = $excel.ActiveWorkbook

$newUrbanBefore1945      = "50% MUR/LWAL+CDN/H:1`n50% MUR/LWAL+CDN/H:2"
$newUrban19461970        = "50% MUR/LWAL+CDN/H:1`n50% MUR/LWAL+CDN/H:2"
$newUrban19712005        = "49% MUR/LWAL+CDN/H:1`n49% MUR/LWAL+CDN/H:2`n1% W/LWAL+CDN/H:1`n1% W/LWAL+CDN/H:2"
$newUrbanAfter2006       = "45% MUR/LWAL+CDN/H:1`n45% MUR/LWAL+CDN/H:2`n5% W/LWAL+CDN/H:1`n5% W/LWAL+CDN/H:2"

foreach ($sheetName in @("mapping_urban", "mapping_rural")) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("F2").Value = $newUrbanBefore1945
    $ws.Range("G2").Value = $newUrban19461970
    $ws.Range("H2").Value = $newUrban19712005
    $ws.Range("I2").Value = $newUrbanAfter2006

    $ws.Range("F2:I2").Font.Size = 11
}

# Make "mapping_urban" the active tab/selection (it was "mapping_rural"
# before), and move the selection to reflect where the edits were made.
$wsUrban = $wb.Worksheets.Item("mapping_urban")
$wsRural = $wb.Worksheets.Item("mapping_rural")

$wsRural.Activate()
$wsRural.Range("F2:I2").Select()

$wsUrban.Activate()
$wsUrban.Range("F4").Select()
